# Use joint centers where possible for k-chain
#
# The "Cylinder Radius" constant used to live on the Poses sheet (J4 label /
# K4 value) and was referenced locally by D2/D4 there. This change moves
# that constant onto the Height sheet (as H3 label / I3 value), folds it
# into the Height sheet's F-column "joint center" k-chain offsets (F5 now
# equals the cylinder radius instead of 0, and F2/F3/F4/F6/F7/F8/F9 all add
# F5 into their chains), and repoints the Poses sheet formulas that used to
# read the local K4 at Height!I3 instead. A couple of Poses-sheet D-column
# formulas that summed Height offsets without the wrist/cylinder radius now
# add D2 (or reference Height!F8 directly) so they include it too.

$wb = $excel.ActiveWorkbook

$wsHeight = $wb.Worksheets.Item("Height")
$wsPoses  = $wb.Worksheets.Item("Poses")

# --- Height sheet: introduce the "Cylinder Radius" constant at H3/I3 ---
$wsHeight.Range("H3").Value = "Cylinder Radius"
$wsHeight.Range("I3").Value = 0.02

# F5 used to be a hardcoded 0; it now picks up the cylinder radius.
$wsHeight.Range("F5").Formula = "=I3"

# Fold F5 (the cylinder radius) into the rest of the k-chain offsets.
$wsHeight.Range("F2").Formula = '=$E$6+$E$7+$E$9+$E$4+$E$3+$E$2/2+F5'
$wsHeight.Range("F3").Formula = '=$E$6+$E$7+$E$9+$E$4 + $E$3/2+F5'
$wsHeight.Range("F4").Formula = '=$E$6+$E$7+$E$9+$E$4/2+F5'
$wsHeight.Range("F6").Formula = "=E6/2+F5"
$wsHeight.Range("F7").Formula = "=E6+E7/2+F5"
$wsHeight.Range("F8").Formula = "=E9+E7+E6+E8/2+F5"
$wsHeight.Range("F9").Formula = "=E6+E7+E9/2+F5"

# --- Poses sheet: point D2/D4 at the relocated constant and drop the old
#     local "Cylinder Radius" label/value (J4/K4) ---
$wsPoses.Range("D2").Formula = "=Height!I3"
$wsPoses.Range("D4").Formula = "=Height!I3"
$wsPoses.Range("J4:K4").ClearContents()

# D14/D18 used to sum Height offsets only up to the shoulder/elbow joint
# centers; now they add back in D2 (the wrist cylinder-radius offset) so the
# hand joint center lines up with the new k-chain.
$wsPoses.Range("D14").Formula = '=Height!$E$6+Height!$E$7+D2'
$wsPoses.Range("D18").Formula = '=Height!$E$6+Height!$E$7+Height!$E$9+D2'

# D20 used to be a hardcoded constant mirroring Height!F8; now it references
# it live so it updates with the rest of the chain.
$wsPoses.Range("D20").Formula = "=Height!F8"

# --- Selections: leave the active cell on each sheet where the author left
#     off, with Poses as the active/selected tab (matches the saved file). ---
$wsHeight.Activate()
$wsHeight.Range("F10").Select()

$wsPoses.Activate()
$wsPoses.Range("D29").Select()

$wb.Save()
